$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the on-call rotation schedule by rolling every start/end date
# forward by one year (365 days), same as the previous run of dates.
$rng = $ws.Range("C3:D28")
foreach ($cell in $rng) {
    $cell.Value = $cell.Value2 + 365
}

$ws.Range("A2").Select()
